$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p106v_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p106v_1</id>", 2)
